# Generate Excel button working: rewrite "Items", "Option Group" and
# "Options" sheets with the new Steak/Larger Plates menu data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Items" — replace rows 2-3 and add rows 4-10
# ---------------------------------------------------------------------
$items = $wb.Worksheets.Item("Items")

$itemRows = @(
    @("Larger Plates", "Nichols Smoked Chicken Breast", "", 36, "Sides"),
    @("Larger Plates", "Scottsdale Pork Belly",          "", 36, "Sides"),
    @("Larger Plates", "Wild Clover Lamb",                "", 40, "Sides"),
    @("Larger Plates", "Rannoch Farm Quail",               "", 37, "Sides"),
    @("Larger Plates", "Pan Roasted Fish",                 "", 37, "Sides"),
    @("Larger Plates", "Braised Lentils",                  "", 30, "Sides"),
    @("Steak Plates",  "300g Porterhouse",                 "", 36, "Sides, Doneness"),
    @("Steak Plates",  "300g Scotch Fillet",                "", 42, "Sides"),
    @("Steak Plates",  "220g Eye Fillet",                    "", 42, "Sides")
)

$r = 2
foreach ($row in $itemRows) {
    $items.Cells.Item($r, 1).Value = $row[0]
    $items.Cells.Item($r, 2).Value = $row[1]
    $items.Cells.Item($r, 3).Value = $row[2]
    $items.Cells.Item($r, 4).Value = $row[3]
    $items.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet "Option Group" — replace rows 2-3
# ---------------------------------------------------------------------
$optionGroup = $wb.Worksheets.Item("Option Group")

$optionGroup.Cells.Item(2, 1).Value = "Doneness"
$optionGroup.Cells.Item(2, 2).Value = "'TRUE"
$optionGroup.Cells.Item(2, 3).Value = "'TRUE"

$optionGroup.Cells.Item(3, 1).Value = "Sides"
$optionGroup.Cells.Item(3, 2).Value = "'FALSE"
$optionGroup.Cells.Item(3, 3).Value = "'TRUE"

# ---------------------------------------------------------------------
# Sheet "Options" — replace rows 2-4 and add rows 5-11
# ---------------------------------------------------------------------
$options = $wb.Worksheets.Item("Options")

$optionRows = @(
    @("Doneness", "Medium Rare", 1),
    @("Doneness", "Medium", 0),
    @("Doneness", "Well Done", 0),
    @("Doneness", "Overcooked", 0),
    @("Doneness", "Burnt", 0),
    @("Sides", "French Fries", 1),
    @("Sides", "Garden Veggie", 1),
    @("Sides", "Coleslaw", 1),
    @("Sides", "Mashed Potato", 1),
    @("Sides", "Tasty Rice", 1)
)

$r = 2
foreach ($row in $optionRows) {
    $options.Cells.Item($r, 1).Value = $row[0]
    $options.Cells.Item($r, 2).Value = $row[1]
    $options.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
